# Update the Handback status report timestamps for the
# e48a681a-322c-41fe-9339-23f16b108803 file across the Overview, zh-cn and
# de-de sheets (regenerated report timestamps).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 4
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-26 20:44:58"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for row 4
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-26 20:44:53"
$wsZhCn.Range("K4").Value = "2016-08-26 20:45:22"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for row 4
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-26 20:44:58"
$wsDeDe.Range("K4").Value = "2016-08-26 20:45:29"
